$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''41.611.53'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.20%  '

$ws.Range('D3').Value = '''2.472.83'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.16%  '

$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').Value = '''318.02'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.53%  '

$ws.Range('D6').Value = '''92.11'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.54%  '

$ws.Range('E7').Value = '  +1.44%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('E9').Value = '  +1.97%  '

$ws.Range('D10').Value = '''0.0863'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +10.07%  '

$ws.Range('D11').Value = '''32.97'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.94%  '

$ws.Range('E12').Value = '  +0.35%  '

$ws.Range('D13').Value = '''2.853.43'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.18%  '

$ws.Range('E14').Value = '  +0.68%  '

$ws.Range('D15').Value = '''15.58'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.53%  '

$ws.Range('D16').Value = '''2.477.85'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.22%  '

$ws.Range('D17').Value = '''0.787'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.45%  '

$ws.Range('D18').Value = '''41.569.24'
$ws.Range('D18').Style = "Normal"

$ws.Range('D19').Value = '''6.49'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.50%  '

$ws.Range('D20').Value = '''0.0₃0949'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.52%  '

$ws.Range('D21').Value = '''71.23'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.42%  '

$ws.Range('D22').Value = '''11.31'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.93%  '

$ws.Range('D23').Value = '''239.66'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.39%  '

$ws.Range('E24').Value = '  +1.20%  '

$ws.Range('E25').Value = '  +1.45%  '

$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').Value = '''24.68'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.31%  '

$ws.Range('E28').Value = '  +3.48%  '

$ws.Range('D29').Value = '''9.86'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.40%  '

$ws.Range('D30').Value = '''36.26'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.27%  '

$ws.Range('D31').Value = '''161.07'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.71%  '

$ws.Range('E32').Value = '  +1.40%  '

$ws.Range('E34').Value = '  +1.65%  '

$ws.Range('E35').Value = '  +0.14%  '

$ws.Range('D36').Value = '''17.23'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.45%  '

$ws.Range('E37').Value = '  +0.36%  '

$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = '''0.116'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.30%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '''1.83'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.75%  '

$ws.Range('E40').Value = '  -3.45%  '

$ws.Range('E41').Value = '  -3.25%  '

$ws.Range('E42').Value = '  +2.66%  '

$ws.Range('D43').Value = '''1.990.66'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.17%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''19.05'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.05%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0285'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.68%  '

$ws.Range('D46').Value = '''2.98'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.29%  '

$ws.Range('D47').Value = '''9.19'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.66%  '

$ws.Range('D48').Value = '''2.711.50'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.23%  '

$ws.Range('D49').Value = '''97.49'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.19%  '

$ws.Range('D50').Value = '''73.96'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.29%  '

$ws.Range('D51').Value = '''66.89'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.53%  '
